$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Ost/Nord (Q/R) coordinate columns for rows 4 and 6 in place.
$ws.Range("Q4").Value = 676709
$ws.Range("R4").Value = 6618511

$ws.Range("Q6").Value = 676487
$ws.Range("R6").Value = 6618440

# The species records that occupied rows 5, 7 and 8 get rotated:
#   row5 <- old row7 data, row7 <- old row8 data, row8 <- old row5 data.
# Capture the original row values before they are overwritten. Use
# .Value2 for the reads -- .Value on this host returns the property
# descriptor text instead of the cell's actual contents.
$a5 = $ws.Range("A5").Value2
$b5 = $ws.Range("B5").Value2
$d5 = $ws.Range("D5").Value2
$e5 = $ws.Range("E5").Value2
$f5 = $ws.Range("F5").Value2
$g5 = $ws.Range("G5").Value2
$h5 = $ws.Range("H5").Value2

$a7 = $ws.Range("A7").Value2
$b7 = $ws.Range("B7").Value2
$d7 = $ws.Range("D7").Value2
$e7 = $ws.Range("E7").Value2
$f7 = $ws.Range("F7").Value2
$g7 = $ws.Range("G7").Value2
$h7 = $ws.Range("H7").Value2

$a8 = $ws.Range("A8").Value2
$b8 = $ws.Range("B8").Value2
$d8 = $ws.Range("D8").Value2
$e8 = $ws.Range("E8").Value2
$f8 = $ws.Range("F8").Value2
$g8 = $ws.Range("G8").Value2
$h8 = $ws.Range("H8").Value2

# Row 5 becomes the old row 7 record.
$ws.Range("A5").Value = $a7
$ws.Range("B5").Value = $b7
$ws.Range("D5").Value = $d7
$ws.Range("E5").Value = $e7
$ws.Range("F5").Value = $f7
$ws.Range("G5").Value = $g7
$ws.Range("H5").Value = $h7
$ws.Range("Q5").Value = 676487
$ws.Range("R5").Value = 6618440

# Row 7 becomes the old row 8 record.
$ws.Range("A7").Value = $a8
$ws.Range("B7").Value = $b8
$ws.Range("D7").Value = $d8
$ws.Range("E7").Value = $e8
$ws.Range("F7").Value = $f8
$ws.Range("G7").Value = $g8
$ws.Range("H7").Value = $h8
$ws.Range("Q7").Value = 676487
$ws.Range("R7").Value = 6618440

# Row 8 becomes the old row 5 record.
$ws.Range("A8").Value = $a5
$ws.Range("B8").Value = $b5
$ws.Range("D8").Value = $d5
$ws.Range("E8").Value = $e5
$ws.Range("F8").Value = $f5
$ws.Range("G8").Value = $g5
$ws.Range("H8").Value = $h5
$ws.Range("Q8").Value = 676709
$ws.Range("R8").Value = 6618511
